$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '92.666.73'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +4.73%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.272.46'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.13%  '

$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '217.86'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.98%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '629.54'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.11%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.401'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.00%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.711'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +3.11%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.999'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.06%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '3.274.14'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.04%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.587'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.43%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000268'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.49%  '

$ws.Range("E13").Value = '  -2.33%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.23'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.60%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '92.406.10'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +4.79%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.879.32'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.04%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.34'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.38%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.265.26'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.05%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.29'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +5.26%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0000221'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +65.34%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.96'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.48%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '447.72'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.63%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.83'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.11%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.25'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -3.05%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.35'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.59%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.10'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.60%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '3.452.84'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.91%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '77.63'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.78%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.00'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.01%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.174'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -4.89%  '

$ws.Range("E31").Value = '  -0.24%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.75'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.66%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '555.80'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.80%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.84'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +28.53%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '7.13'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.22%  '

$ws.Range("E36").Value = '  -2.06%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.28'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -8.72%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '22.62'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.27%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '22.46'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +3.12%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.130'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -6.42%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.00'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.04%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.393'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.33%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.99'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.84%  '

$ws.Range("E44").Value = '  -0.08%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '150.26'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.66%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '45.34'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.20%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '178.99'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.01%  '

$ws.Range("E48").Value = '  +2.10%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.28'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.33%  '

$ws.Range("E50").Value = '  +1.87%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '4.22'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.64%  '
